{"js": "// Replace each two-digit multiplication equation in the table with its\n// updated counterpart. Every \"old\" string occurs exactly once in the\n// document, so a targeted search + replace is used per pair (rather than\n// a blind sequential replace) to stay robust to ordering.\nconst replacements = [\n  [\"46\u00d711=506\", \"32\u00d726=832\"],\n  [\"93\u00d745=4185\", \"22\u00d755=1210\"],\n  [\"27\u00d782=2214\", \"51\u00d735=1785\"],\n  [\"35\u00d775=2625\", \"15\u00d752=780\"],\n  [\"87\u00d763=5481\", \"96\u00d731=2976\"],\n  [\"19\u00d778=1482\", \"99\u00d786=8514\"],\n  [\"26\u00d769=1794\", \"13\u00d738=494\"],\n  [\"17\u00d784=1428\", \"67\u00d748=3216\"],\n  [\"12\u00d744=528\", \"50\u00d785=4250\"],\n  [\"71\u00d762=4402\", \"63\u00d721=1323\"],\n  [\"45\u00d732=1440\", \"25\u00d764=1600\"],\n  [\"74\u00d713=962\", \"83\u00d733=2739\"],\n  [\"88\u00d772=6336\", \"50\u00d751=2550\"],\n  [\"29\u00d757=1653\", \"67\u00d720=1340\"],\n  [\"18\u00d752=936\", \"11\u00d752=572\"],\n  [\"97\u00d750=4850\", \"45\u00d728=1260\"],\n  [\"96\u00d798=9408\", \"26\u00d729=754\"],\n  [\"25\u00d741=1025\", \"39\u00d751=1989\"],\n  [\"81\u00d744=3564\", \"14\u00d747=658\"],\n  [\"75\u00d779=5925\", \"92\u00d760=5520\"],\n  [\"37\u00d749=1813\", \"79\u00d778=6162\"],\n  [\"35\u00d762=2170\", \"81\u00d771=5751\"],\n  [\"16\u00d785=1360\", \"31\u00d716=496\"],\n  [\"39\u00d749=1911\", \"30\u00d768=2040\"],\n  [\"69\u00d727=1863\", \"87\u00d769=6003\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication equation in the table with its\n# updated counterpart. Each \"old\" equation string occurs exactly once in\n# the document, so Find/Replace (ReplaceAll) per pair is safe and precise.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"46\u00d711=506\", \"32\u00d726=832\"),\n    @(\"93\u00d745=4185\", \"22\u00d755=1210\"),\n    @(\"27\u00d782=2214\", \"51\u00d735=1785\"),\n    @(\"35\u00d775=2625\", \"15\u00d752=780\"),\n    @(\"87\u00d763=5481\", \"96\u00d731=2976\"),\n    @(\"19\u00d778=1482\", \"99\u00d786=8514\"),\n    @(\"26\u00d769=1794\", \"13\u00d738=494\"),\n    @(\"17\u00d784=1428\", \"67\u00d748=3216\"),\n    @(\"12\u00d744=528\", \"50\u00d785=4250\"),\n    @(\"71\u00d762=4402\", \"63\u00d721=1323\"),\n    @(\"45\u00d732=1440\", \"25\u00d764=1600\"),\n    @(\"74\u00d713=962\", \"83\u00d733=2739\"),\n    @(\"88\u00d772=6336\", \"50\u00d751=2550\"),\n    @(\"29\u00d757=1653\", \"67\u00d720=1340\"),\n    @(\"18\u00d752=936\", \"11\u00d752=572\"),\n    @(\"97\u00d750=4850\", \"45\u00d728=1260\"),\n    @(\"96\u00d798=9408\", \"26\u00d729=754\"),\n    @(\"25\u00d741=1025\", \"39\u00d751=1989\"),\n    @(\"81\u00d744=3564\", \"14\u00d747=658\"),\n    @(\"75\u00d779=5925\", \"92\u00d760=5520\"),\n    @(\"37\u00d749=1813\", \"79\u00d778=6162\"),\n    @(\"35\u00d762=2170\", \"81\u00d771=5751\"),\n    @(\"16\u00d785=1360\", \"31\u00d716=496\"),\n    @(\"39\u00d749=1911\", \"30\u00d768=2040\"),\n    @(\"69\u00d727=1863\", \"87\u00d769=6003\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
